$d = $word.ActiveDocument

# --- "Dean:" paragraph: the dean's name changes from
#     "Caryn L. Beck-Dudley" to "Ed Grier". The space that used to sit in
#     front of the name now trails the "Dean:" label instead, so the
#     rendered line still reads "Dean: Ed Grier". ---
$deanPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Dean:")) {
        $deanPara = $p
        break
    }
}

if ($deanPara -ne $null) {
    $paraStart = $deanPara.Range.Start
    $paraEnd = $deanPara.Range.End
    $paraText = $deanPara.Range.Text
    $splitOffset = $paraText.IndexOf(" Caryn L. Beck-Dudley")
    $labelEnd = $paraStart + $splitOffset

    # Replace the name (tail) before the label (head) so the head's
    # offsets -- computed above -- are still valid when used.
    $nameRange = $d.Range($labelEnd, $paraEnd - 1)
    $nameRange.Text = "Ed Grier "

    $labelRange = $d.Range($paraStart, $labelEnd)
    $labelRange.Text = "Dean: "
}

# --- Explicitly mark the page section as portrait orientation ---
$d.PageSetup.Orientation = 0
